$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts existing rows 7-14 down to 8-15)
# to hold the new "populationsCSV" property.
$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = "populationsCSV"
$ws.Range("B7").Value = "PopulationsCSV"
$ws.Range("C7").Value = "Name of the folder containing population defined in files"

# Update the active selection to match the saved view state.
$ws.Range("G30").Select()
